# shopping_mall_page_info.xlsx
# - add NAVER Shopping (네이버쇼핑) information as a new data row
# - add a "parent element Xpath" (부모 노드 Xpath) column, shifting the
#   existing product-name / price Xpath columns one column to the right
#
# Shared-string order matters for byte-faithful output: the statements
# below are ordered so brand-new text values are written to the
# worksheet in the same order the original author's Excel session
# would have produced them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Header row (row 2): move "상품명 Xpath" one column right (I2),
#    introduce the new "부모 노드 Xpath" header (H2), and restore the
#    "가격 Xpath" header in the new last column (J2).
# ----------------------------------------------------------------------
$ws.Range("I2").Value = "상품명 Xpath"
$ws.Range("H2").Value = "부모 노드 Xpath"
$ws.Range("J2").Value = "가격 Xpath"

# ----------------------------------------------------------------------
# 2) Tmon (티몬) data row (row 3): switch the product/price Xpaths to
#    be relative to the new parent-node Xpath, and record that parent
#    node Xpath itself.
# ----------------------------------------------------------------------
$ws.Range("I3").Value = "./p[@class='title']/strong[@class='tx']"
$ws.Range("J3").Value = "./div[@class='price_area']/span[@class='price']/span[@class='sale']/i[@class='num']"
$ws.Range("H3").Value = "//*[@id='search_app']/div[2]/section/div/ul/div/div/li/a/div[3]"

# ----------------------------------------------------------------------
# 3) New NAVER Shopping (네이버쇼핑) data row (row 4).
# ----------------------------------------------------------------------
$ws.Range("C4").Value = "네이버쇼핑"
$ws.Range("D4").Value = "https://search.shopping.naver.com/search/"
$ws.Range("E4").Value = "all.nhn?query={검색어}"
$ws.Range("F4").Value = "사용 불가"
$ws.Range("H4").Value = "//*[@id=`"_search_list`"]/div[1]/ul/li/div[@class='info']/"
$ws.Range("I4").Value = "./div[@class='tit']/a"
$ws.Range("J4").Value = "./span[@class='price']/em/span[@class='num _price_reload'] OR, span[@class='price']/em/span[@class='num']"
$ws.Range("G4").Value = "[todo] 추가하기"

# ----------------------------------------------------------------------
# 4) Formatting - reuse the existing style slots (header / hyperlink /
#    "price Xpath" emphasis) instead of letting Excel fabricate new
#    cell-format entries, by copy/pasting formats from cells that
#    already carry the right style.
# ----------------------------------------------------------------------

# Header style (gray fill, centered) onto the new header cells.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("H2:J2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The "price Xpath" emphasis style moves from I3 to J3 (must happen
# before I3's own style gets reset below).
$ws.Range("I3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Plain/default style onto the Tmon row's relocated Xpath cells
# (they must lose whatever style they inherited from their old
# column position).
$ws.Range("C3").Copy() | Out-Null
$ws.Range("H3:I3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Plain/default style for the whole new NAVER row ...
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4:J4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ... except the search-URL cell, which follows the same hyperlink
# style as D3.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 5) Hyperlink for the new NAVER Shopping search URL.
# ----------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D4"), "https://search.shopping.naver.com/search/") | Out-Null

# Re-apply the hyperlink-style formatting (Hyperlinks.Add resets it to
# a fresh style slot) so D4 keeps sharing D3's existing style.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 6) Column widths (best-fit widths recorded by Excel for the wider
#    new content).
# ----------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 38.0
$ws.Columns.Item(5).ColumnWidth = 20.571428571428573
$ws.Columns.Item(8).ColumnWidth = 57.714285714285715
$ws.Columns.Item(9).ColumnWidth = 54.57142857142857
$ws.Columns.Item(10).ColumnWidth = 98.28571428571429

# ----------------------------------------------------------------------
# 7) View state: scroll the window so column H is left-most, and leave
#    the active selection at J15 like the original author's session.
# ----------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("J15").Select() | Out-Null
